# step-showcase.xlsx update
# - [json] add new function storeKeys(json,jsonpath,var), inserted in its
#   alphabetically-sorted slot (between storeCount and storeValue) on the
#   hidden "#system" lookup sheet.
# - remove the (now orphaned) "text" category column from the "#system"
#   sheet and its entry from the "target" category list; every category
#   column to the right of it (web, webalert, webcookie, ws, ws.async,
#   xml) shifts one column to the left.
# - keep every affected named range ("json", "target", "web", "webalert",
#   "webcookie", "ws", "ws.async", "xml") in sync with its new address.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) [json] insert storeKeys(json,jsonpath,var) at M16, pushing the two
#    rows below it (storeValue / storeValues) down by one row. We avoid
#    Range.Insert/Delete here since (in this host) they shift the whole
#    row rather than just the target column - so the move is done with
#    plain cell-value assignment instead, touching column M only.
# ---------------------------------------------------------------------
$oldM16 = $ws.Range("M16").Value()
$oldM17 = $ws.Range("M17").Value()
$ws.Range("M18").Value = $oldM17
$ws.Range("M17").Value = $oldM16
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 2) Drop the "text" category column (Y) outright - this is a genuine
#    full-column delete so every later category column (web=Z,
#    webalert=AA, webcookie=AB, ws=AC, ws.async=AD, xml=AE) shifts left
#    by one (-> Y, Z, AA, AB, AC, AD respectively).
# ---------------------------------------------------------------------
$ws.Columns("Y").Delete()

# ---------------------------------------------------------------------
# 3) Remove the now-stale "text" row from the "target" category index
#    in column A (was A25), shifting web..xml (A26:A31) up by one row
#    and clearing the now-unused last row (A31). Again done with plain
#    value assignment so only column A is touched.
# ---------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $below = $ws.Cells.Item($r + 1, 1).Value()
    $ws.Cells.Item($r, 1).Value = $below
}
$ws.Cells.Item(31, 1).ClearContents()

# ---------------------------------------------------------------------
# 4) Re-point the named ranges that moved/resized as a result of the
#    edits above.
# ---------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
